$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 438-440, pushing the existing rows 438-542 down to 441-545.
$ws.Rows("438:440").Insert()

# Common (constant) columns shared by every data row in this table.
$colA = 3
$colB = "Femacal de La Calera"
$colC = "Coquimbo"
$colE = 5
$colF = "Fruta"
$colG = 100102
$colH = "Cítricos"
$colI = 100102004
$colJ = "Mandarina"
$colQ = "`$/bandeja 10 kilos"
$colR = "Provincia de Quillota"
$colT = 10

# New row 438: Murcott / Especial
$ws.Range("A438").Value = $colA
$ws.Range("B438").Value = $colB
$ws.Range("C438").Value = $colC
$ws.Range("D438").Value = 44508
$ws.Range("E438").Value = $colE
$ws.Range("F438").Value = $colF
$ws.Range("G438").Value = $colG
$ws.Range("H438").Value = $colH
$ws.Range("I438").Value = $colI
$ws.Range("J438").Value = $colJ
$ws.Range("K438").Value = "Murcott"
$ws.Range("L438").Value = "Especial"
$ws.Range("M438").Value = 60
$ws.Range("N438").Value = 6000
$ws.Range("O438").Value = 6000
$ws.Range("P438").Value = 6000
$ws.Range("Q438").Value = $colQ
$ws.Range("R438").Value = $colR
$ws.Range("S438").Value = 600
$ws.Range("T438").Value = $colT

# New row 439: Murcott / Primera
$ws.Range("A439").Value = $colA
$ws.Range("B439").Value = $colB
$ws.Range("C439").Value = $colC
$ws.Range("D439").Value = 44508
$ws.Range("E439").Value = $colE
$ws.Range("F439").Value = $colF
$ws.Range("G439").Value = $colG
$ws.Range("H439").Value = $colH
$ws.Range("I439").Value = $colI
$ws.Range("J439").Value = $colJ
$ws.Range("K439").Value = "Murcott"
$ws.Range("L439").Value = "Primera"
$ws.Range("M439").Value = 67
$ws.Range("N439").Value = 5000
$ws.Range("O439").Value = 5000
$ws.Range("P439").Value = 5000
$ws.Range("Q439").Value = $colQ
$ws.Range("R439").Value = $colR
$ws.Range("S439").Value = 500
$ws.Range("T439").Value = $colT

# New row 440: Murcott / Segunda
$ws.Range("A440").Value = $colA
$ws.Range("B440").Value = $colB
$ws.Range("C440").Value = $colC
$ws.Range("D440").Value = 44508
$ws.Range("E440").Value = $colE
$ws.Range("F440").Value = $colF
$ws.Range("G440").Value = $colG
$ws.Range("H440").Value = $colH
$ws.Range("I440").Value = $colI
$ws.Range("J440").Value = $colJ
$ws.Range("K440").Value = "Murcott"
$ws.Range("L440").Value = "Segunda"
$ws.Range("M440").Value = 65
$ws.Range("N440").Value = 4000
$ws.Range("O440").Value = 4000
$ws.Range("P440").Value = 4000
$ws.Range("Q440").Value = $colQ
$ws.Range("R440").Value = $colR
$ws.Range("S440").Value = 400
$ws.Range("T440").Value = $colT
